# Generate Report for Archive
#
# 1) Replace the "Ready for handoff" status text with "In Translation"
#    everywhere it appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2) Narrow the "Status" / language columns that held that text from
#    ~17.22 chars wide to ~13.41 chars wide (Overview columns E & F,
#    and column C on the zh-cn and de-de sheets).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newWidth = 13.4101845877511

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value2 -eq $oldStatus) {
                $cell.Value2 = $newStatus
            }
        }
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newWidth
